$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the MSIDN (col B), SERIAL (col C) and PLU (col D) values for the last
# two data rows, column by column (matches the order the new values were
# appended to the shared-string table).
$ws.Cells.Item(11, 2).Value = "8957732111198172291"
$ws.Cells.Item(12, 2).Value = "8957732111198172290"

$ws.Cells.Item(11, 3).Value = "3016875982"
$ws.Cells.Item(12, 3).Value = "3016875893"

$ws.Cells.Item(11, 4).Value = "732111198172291"
$ws.Cells.Item(12, 4).Value = "732111198172290"

# Remove row 13 entirely (its MSIDN/SERIAL/PLU entry is dropped from the sheet).
$ws.Rows.Item(13).Delete()

# Reflect the updated active selection left by the edit.
$ws.Range("B13").Select()
